# Updates scraped betexplorer odds data for serbia/prva-liga 2023-2024:
#  - rows 111-113 (11/11/2023 matches) reshuffled to correct order
#  - rows 119-121 (18/11/2023 matches) reshuffled to correct order
#  - rows 137-138 (27/11/2023 matches) swapped
#  - new row 139 appended for Kolubara vs Sloboda (27/11/2023)

function Set-RowFV {
    param($ws, $row, $values)
    $col = 6
    foreach ($v in $values) {
        $ws.Cells.Item($row, $col).Value = $v
        $col = $col + 1
    }
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Group 1: 11/11/2023 matches, rows 111-113 rotated ---
$row111 = @("Radnicki S. Mitrovica", 0, "Mladost GAT", 0, 2.56, "11/11/2023 02:13", 2.04, "11/11/2023 12:53", 2.74, "11/11/2023 02:13", 2.73, "11/11/2023 12:53", 2.81, "11/11/2023 02:13", 4.03, "11/11/2023 12:53", "https://www.betexplorer.com/football/serbia/prva-liga/radnicki-s-mitrovica-mladost-gat/88pX5ob4/")
$row112 = @("Tekstilac Odzaci", 1, "Macva", 0, 1.8, "11/11/2023 02:13", 2.15, "11/11/2023 12:58", 2.98, "11/11/2023 02:13", 2.41, "11/11/2023 12:58", 4.32, "11/11/2023 02:13", 4.44, "11/11/2023 12:58", "https://www.betexplorer.com/football/serbia/prva-liga/tekstilac-odzaci-macva-sabac/2Vqy55DA/")
$row113 = @("FK Indjija", 1, "Smederevo", 2, 1.88, "11/11/2023 02:13", 1.99, "11/11/2023 12:48", 2.98, "11/11/2023 02:13", 2.8, "11/11/2023 12:48", 3.91, "11/11/2023 02:13", 4.05, "11/11/2023 12:35", "https://www.betexplorer.com/football/serbia/prva-liga/indjija-smederevo/UHoT6Rrb/")

Set-RowFV $ws 111 $row111
Set-RowFV $ws 112 $row112
Set-RowFV $ws 113 $row113

# --- Group 2: 18/11/2023 matches, rows 119-121 rotated ---
$row119 = @("Metalac", 2, "OFK Beograd", 0, 4.26, "18/11/2023 02:12", 4.21, "18/11/2023 12:55", 3.19, "18/11/2023 02:12", 3.32, "18/11/2023 12:55", 1.74, "18/11/2023 02:12", 1.76, "18/11/2023 12:55", "https://www.betexplorer.com/football/serbia/prva-liga/metalac-ofk-beograd/jkoC6IlF/")
$row120 = @("Jedinstvo U.", 3, "Graficar Beograd", 1, 1.98, "18/11/2023 02:12", 1.95, "18/11/2023 12:51", 3.07, "18/11/2023 02:12", 3.27, "18/11/2023 12:51", 3.45, "18/11/2023 02:12", 3.49, "18/11/2023 12:51", "https://www.betexplorer.com/football/serbia/prva-liga/jedinstvo-ub-graficar-beograd/dSo87bZ8/")
$row121 = @("Vrsac", 1, "Radnicki S. Mitrovica", 0, 2.21, "18/11/2023 02:12", 2.43, "18/11/2023 12:30", 2.78, "18/11/2023 02:12", 2.64, "18/11/2023 12:30", 3.25, "18/11/2023 02:12", 3.18, "18/11/2023 12:30", "https://www.betexplorer.com/football/serbia/prva-liga/vrsac-radnicki-s-mitrovica/Ai1D4dJR/")

Set-RowFV $ws 119 $row119
Set-RowFV $ws 120 $row120
Set-RowFV $ws 121 $row121

# --- Group 3: 27/11/2023 matches, rows 137-138 swapped ---
$row137 = @("RFK Novi Sad", 0, "Dubocica", 0, 2.82, "26/11/2023 02:13", 2.83, "27/11/2023 12:55", 2.77, "26/11/2023 02:13", 2.75, "27/11/2023 12:55", 2.47, "26/11/2023 02:13", 2.58, "27/11/2023 12:55", "https://www.betexplorer.com/football/serbia/prva-liga/rfk-novi-sad-dubocica/jBMnJeYK/")
$row138 = @("Radnicki S. Mitrovica", 1, "Macva", 0, 2.05, "27/11/2023 01:12", 2.11, "27/11/2023 12:51", 2.78, "27/11/2023 01:12", 2.63, "27/11/2023 12:51", 3.66, "27/11/2023 01:12", 3.99, "27/11/2023 12:01", "https://www.betexplorer.com/football/serbia/prva-liga/radnicki-s-mitrovica-macva-sabac/2kthcEB1/")

Set-RowFV $ws 137 $row137
Set-RowFV $ws 138 $row138

# --- New row 139: Kolubara vs Sloboda, 27/11/2023 18:00 ---
# Copy formatting from the row above so styles (bold/border/centered index, date format) match.
$ws.Range("A138").Copy()
$ws.Range("A139").PasteSpecial(-4122)
$ws.Range("E138").Copy()
$ws.Range("E139").PasteSpecial(-4122)

$ws.Cells.Item(139, 1).Value = 138
$ws.Cells.Item(139, 2).Value = "serbia"
$ws.Cells.Item(139, 3).Value = "prva-liga"
$ws.Cells.Item(139, 4).Value = "2023-2024"
$ws.Cells.Item(139, 5).Value = 45257.75

$row139 = @("Kolubara", 0, "Sloboda", 0, 1.57, "27/11/2023 07:12", 1.54, "27/11/2023 17:54", 3.33, "27/11/2023 07:12", 3.52, "27/11/2023 17:54", 5.37, "27/11/2023 07:12", 5.76, "27/11/2023 17:54", "https://www.betexplorer.com/football/serbia/prva-liga/kolubara-sloboda/GdQwQG1B/")
Set-RowFV $ws 139 $row139
